# Updates the cryptos price/volume table (and two name/link row swaps)
# to reflect the latest scraped values, keeping every cell as plain text
# (matching the original inlineStr cells) instead of letting Excel
# auto-convert numeric- or percent-looking strings into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Belt-and-braces: stop Excel from auto-converting "12.3%"-like input
# into a real percentage value when we assign it below.
$excel.AutoPercentEntry = $false

$updates = @(
    @(2, 4, "58.305.65"),
    @(2, 5, "  -3.15%  "),
    @(3, 4, "2.450.17"),
    @(3, 5, "  -4.04%  "),
    @(4, 5, "  +0.01%  "),
    @(5, 4, "528.42"),
    @(6, 4, "133.65"),
    @(6, 5, "  -7.50%  "),
    @(7, 5, "  +0.30%  "),
    @(8, 4, "0.555"),
    @(8, 5, "  -2.89%  "),
    @(9, 4, "2.454.91"),
    @(9, 5, "  -4.55%  "),
    @(10, 4, "0.0988"),
    @(10, 5, "  -2.34%  "),
    @(11, 5, "  -0.44%  "),
    @(12, 4, "5.29"),
    @(12, 5, "  -3.58%  "),
    @(13, 5, "  -5.60%  "),
    @(14, 4, "2.886.35"),
    @(14, 5, "  -3.84%  "),
    @(15, 4, "58.254.79"),
    @(15, 5, "  -3.17%  "),
    @(16, 5, "  -6.05%  "),
    @(17, 4, "0.0000138"),
    @(17, 5, "  -3.85%  "),
    @(18, 4, "2.460.69"),
    @(18, 5, "  -3.77%  "),
    @(19, 4, "10.74"),
    @(19, 5, "  -4.53%  "),
    @(20, 4, "4.19"),
    @(20, 5, "  -3.19%  "),
    @(21, 4, "320.34"),
    @(21, 5, "  -2.04%  "),
    @(22, 4, "0.996"),
    @(22, 5, "  -0.30%  "),
    @(23, 5, "  -4.45%  "),
    @(24, 4, "62.39"),
    @(24, 5, "  -1.80%  "),
    @(25, 4, "0.407"),
    @(25, 5, "  -6.12%  "),
    @(26, 4, "0.163"),
    @(26, 5, "  -2.30%  "),
    @(27, 5, "  -0.98%  "),
    @(28, 4, "7.43"),
    @(28, 5, "  -7.53%  "),
    @(29, 4, "0.0₃0749"),
    @(29, 5, "  -5.82%  "),
    @(30, 4, "6.49"),
    @(30, 5, "  -8.10%  "),
    @(31, 4, "1.74"),
    @(31, 5, "  -4.00%  "),
    @(32, 4, "164.08"),
    @(32, 5, "  -1.26%  "),
    @(33, 4, "0.999"),
    @(33, 5, "  +0.06%  "),
    @(34, 4, "1.11"),
    @(34, 5, "  -6.62%  "),
    @(35, 2, "EthereumClassic"),
    @(35, 3, "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @(35, 4, "18.19"),
    @(35, 5, "  -2.98%  "),
    @(36, 2, "ImmutableX"),
    @(36, 3, "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @(36, 4, "1.35"),
    @(36, 5, "  -8.83%  "),
    @(37, 5, "  -8.92%  "),
    @(38, 5, "  -6.48%  "),
    @(39, 4, "36.40"),
    @(39, 5, "  -1.79%  "),
    @(40, 4, "0.803"),
    @(40, 5, "  -3.69%  "),
    @(41, 5, "  -5.11%  "),
    @(42, 2, "FirstDigitalUSD"),
    @(42, 3, "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"),
    @(42, 4, "0.997"),
    @(42, 5, "  +0.38%  "),
    @(43, 4, "273.21"),
    @(43, 5, "  -9.33%  "),
    @(44, 2, "RenderToken"),
    @(44, 3, "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @(44, 4, "5.05"),
    @(44, 5, "  -9.99%  "),
    @(45, 4, "10.84"),
    @(45, 5, "  -0.17%  "),
    @(46, 4, "0.584"),
    @(46, 5, "  -4.56%  "),
    @(47, 4, "0.0921"),
    @(47, 5, "  -1.93%  "),
    @(48, 4, "120.70"),
    @(48, 5, "  -5.41%  "),
    @(49, 4, "0.0503"),
    @(49, 5, "  -2.94%  "),
    @(50, 5, "  -5.33%  "),
    @(51, 4, "17.05"),
    @(51, 5, "  -6.70%  "),
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so strings like "528.42" or "18.19" are not
    # silently turned into numeric values, then restore the default "Normal"
    # style so the cell's style index matches the untouched cells around it.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
